$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-CellText 2 4 "62.073.92"
Set-CellText 2 5 "  -7.00%  "
Set-CellText 3 4 "2.890.31"
Set-CellText 3 5 "  -10.13%  "
Set-CellText 4 4 "0.997"
Set-CellText 4 5 "  -0.27%  "
Set-CellText 5 4 "527.40"
Set-CellText 5 5 "  -11.58%  "
Set-CellText 6 4 "127.47"
Set-CellText 6 5 "  -16.25%  "
Set-CellText 7 4 "0.999"
Set-CellText 7 5 "  -0.05%  "
Set-CellText 8 4 "2.869.18"
Set-CellText 8 5 "  -10.54%  "
Set-CellText 9 4 "0.450"
Set-CellText 9 5 "  -17.80%  "
Set-CellText 10 4 "0.138"
Set-CellText 10 5 "  -20.71%  "
Set-CellText 11 4 "5.70"
Set-CellText 11 5 "  -14.46%  "
Set-CellText 12 4 "0.417"
Set-CellText 12 5 "  -16.43%  "
Set-CellText 13 4 "30.74"
Set-CellText 13 5 "  -21.65%  "
Set-CellText 14 4 "0.0000192"
Set-CellText 14 5 "  -21.43%  "
Set-CellText 15 4 "3.368.69"
Set-CellText 15 5 "  -9.89%  "
Set-CellText 16 4 "61.914.88"
Set-CellText 16 5 "  -7.31%  "
Set-CellText 17 5 "  -5.64%  "
Set-CellText 18 4 "2.886.95"
Set-CellText 18 5 "  -10.33%  "
Set-CellText 19 4 "459.36"
Set-CellText 19 5 "  -14.01%  "
Set-CellText 20 4 "5.99"
Set-CellText 20 5 "  -16.67%  "
Set-CellText 21 4 "12.27"
Set-CellText 21 5 "  -18.36%  "
Set-CellText 22 4 "0.608"
Set-CellText 22 5 "  -20.44%  "
Set-CellText 23 2 "Uniswap"
Set-CellText 23 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-CellText 23 4 "6.14"
Set-CellText 23 5 "  -22.79%  "
Set-CellText 24 2 "Litecoin"
Set-CellText 24 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-CellText 24 4 "72.99"
Set-CellText 24 5 "  -15.36%  "
Set-CellText 25 2 "Dai"
Set-CellText 25 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText 25 4 "1.00"
Set-CellText 25 5 "  +0.01%  "
Set-CellText 26 2 "InternetComputer(DFINITY)"
Set-CellText 26 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText 26 4 "11.55"
Set-CellText 26 5 "  -16.92%  "
Set-CellText 27 4 "2.47"
Set-CellText 27 5 "  -23.10%  "
Set-CellText 28 4 "6.68"
Set-CellText 28 5 "  -18.46%  "
Set-CellText 29 4 "1.78"
Set-CellText 29 5 "  -19.07%  "
Set-CellText 30 4 "23.66"
Set-CellText 30 5 "  -19.77%  "
Set-CellText 31 4 "0.996"
Set-CellText 31 5 "  -0.73%  "
Set-CellText 32 5 "  -11.74%  "
Set-CellText 33 4 "2.23"
Set-CellText 33 5 "  -16.91%  "
Set-CellText 34 4 "50.23"
Set-CellText 34 5 "  -5.96%  "
Set-CellText 35 4 "453.52"
Set-CellText 35 5 "  -17.20%  "
Set-CellText 36 4 "5.21"
Set-CellText 36 5 "  -20.27%  "
Set-CellText 37 4 "4.46"
Set-CellText 37 5 "  -22.24%  "
Set-CellText 38 4 "0.0369"
Set-CellText 38 5 "  -14.17%  "
Set-CellText 39 2 "Kaspa"
Set-CellText 39 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText 39 4 "0.110"
Set-CellText 39 5 "  -12.27%  "
Set-CellText 40 2 "Hedera"
Set-CellText 40 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText 40 4 "0.0715"
Set-CellText 40 5 "  -17.07%  "
Set-CellText 41 4 "7.45"
Set-CellText 41 5 "  -20.33%  "
Set-CellText 42 4 "2.559.52"
Set-CellText 42 5 "  -12.51%  "
Set-CellText 43 5 "  -0.26%  "
Set-CellText 44 4 "2.06"
Set-CellText 44 5 "  -24.04%  "
Set-CellText 45 2 "TheGraph"
Set-CellText 45 3 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-CellText 45 4 "0.209"
Set-CellText 45 5 "  -21.01%  "
Set-CellText 46 2 "Stellar"
Set-CellText 46 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-CellText 46 4 "0.0971"
Set-CellText 46 5 "  -15.19%  "
Set-CellText 47 2 "Monero"
Set-CellText 47 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText 47 4 "106.40"
Set-CellText 47 5 "  -11.47%  "
Set-CellText 48 2 "Fetch.AI"
Set-CellText 48 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText 48 4 "1.71"
Set-CellText 48 5 "  -20.69%  "
Set-CellText 49 2 "BitgetToken"
Set-CellText 49 3 "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
Set-CellText 49 4 "1.16"
Set-CellText 49 5 "  -6.51%  "
Set-CellText 50 4 "0.0₃0441"
Set-CellText 50 5 "  -24.79%  "
Set-CellText 51 4 "20.37"
Set-CellText 51 5 "  -23.51%  "
